$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the sampling frequency value
$ws.Range("C1").Value = 250

# New BCI data file timestamps
$ws.Range("F1").Value = 5908
$ws.Range("F2").Value = 5719
$ws.Range("F3").Formula = "=F1-F2"

# C4 becomes a formula derived from the new interval instead of a hard-coded literal
$ws.Range("C4").Formula = "=F3"

# Update the selected cell to reflect where the user ended up
$ws.Range("C5").Select()
